# "Generate Report for Handback"
# Adds a freshly-handed-back file (c70769f7-ccce-4712-8540-acb6e1f4b28a.md) to the
# handback-status report, alongside refreshed timestamps/filenames for the file
# that was already in the report (7df779dc-... -> 0960d0ef-...).

$wb = $excel.ActiveWorkbook

$oldGuid = "7df779dc-91c1-4232-8744-4f4b24bdf9ff"
$guid1   = "0960d0ef-5e52-480e-909c-709439d89262"
$guid2   = "c70769f7-ccce-4712-8540-acb6e1f4b28a"

$zhcnXlf1 = "$guid1.f11c47888053a987de3bb45707e4c80e2fbd9dd9.zh-cn.xlf"
$zhcnXlf2 = "$guid2.aa355650280fec022eaba1aa00bc9145b946d3ec.zh-cn.xlf"
$dedeXlf1 = "$guid1.f11c47888053a987de3bb45707e4c80e2fbd9dd9.de-de.xlf"
$dedeXlf2 = "$guid2.aa355650280fec022eaba1aa00bc9145b946d3ec.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")
$loO = $wsO.ListObjects.Item(1)

# refresh the existing row for guid1 (renamed file + later handback time)
$wsO.Range("B2").Hyperlinks.Delete()
$wsO.Range("A2").Value = "$guid1.md"
$wsO.Range("B2").Value = "e2e\$guid1.md"
$wsO.Range("G2").Value = "2016-11-15 17:37:27"
$wsO.Hyperlinks.Add($wsO.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$guid1.md", "", "", "e2e\$guid1.md")

# append the new row for guid2
$loO.ListRows.Add() | Out-Null
$wsO.Range("A3").Value = "$guid2.md"
$wsO.Range("B3").Value = "e2e\$guid2.md"
$wsO.Range("C3").Value = ".md"
$wsO.Range("E3").Value = "Handed back: in sync with en-US"
$wsO.Range("F3").Value = "Handed back: in sync with en-US"
$wsO.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsO.Range("G3").Value = "2016-11-15 17:38:26"
$wsO.Hyperlinks.Add($wsO.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$guid2.md", "", "", "e2e\$guid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZ = $wb.Worksheets.Item("zh-cn")
$loZ = $wsZ.ListObjects.Item(1)

$wsZ.Range("A2").Hyperlinks.Delete()
$wsZ.Range("I2").Hyperlinks.Delete()
$wsZ.Range("A2").Value = "$guid1.md"
$wsZ.Range("G2").Value = $zhcnXlf1
$wsZ.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("H2").Value = "2016-11-15 17:37:13"
$wsZ.Range("I2").Value = "$guid1.md"
$wsZ.Range("J2").Value = $zhcnXlf1
$wsZ.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("K2").Value = "2016-11-15 17:38:59"
$wsZ.Hyperlinks.Add($wsZ.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$guid1.md", "", "", "$guid1.md")
$wsZ.Hyperlinks.Add($wsZ.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f58de31fe319cf0e756137a19264bec524b8c064/e2e/$guid1.md", "", "", "$guid1.md")

$loZ.ListRows.Add() | Out-Null
$wsZ.Range("B3").Value = ".md"
$wsZ.Range("C3").Value = "Handed back: in sync with en-US"
$wsZ.Range("D3").Value = "e2e"
$wsZ.Range("E3").Value = "ht"
$wsZ.Range("F3").Value = "'True"
$wsZ.Range("G3").Value = $zhcnXlf2
$wsZ.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("H3").Value = "2016-11-15 17:38:12"
$wsZ.Range("J3").Value = $zhcnXlf2
$wsZ.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZ.Range("K3").Value = "2016-11-15 17:38:59"
$wsZ.Range("L3").Value = "'"
$wsZ.Range("M3").Value = "'True"
$wsZ.Range("N3").Value = "'"
$wsZ.Range("O3").Value = "'True"
$wsZ.Range("P3").Value = "'"
$wsZ.Range("A3").Value = "$guid2.md"
$wsZ.Range("I3").Value = "$guid2.md"
$wsZ.Hyperlinks.Add($wsZ.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$guid2.md", "", "", "$guid2.md")
$wsZ.Hyperlinks.Add($wsZ.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f58de31fe319cf0e756137a19264bec524b8c064/e2e/$guid2.md", "", "", "$guid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("de-de")
$loD = $wsD.ListObjects.Item(1)

$wsD.Range("A2").Hyperlinks.Delete()
$wsD.Range("I2").Hyperlinks.Delete()
$wsD.Range("A2").Value = "$guid1.md"
$wsD.Range("G2").Value = $dedeXlf1
$wsD.Range("H2").Value = "2016-11-15 17:37:27"
$wsD.Range("I2").Value = "$guid1.md"
$wsD.Range("J2").Value = $dedeXlf1
$wsD.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("K2").Value = "2016-11-15 17:39:18"
$wsD.Hyperlinks.Add($wsD.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$guid1.md", "", "", "$guid1.md")
$wsD.Hyperlinks.Add($wsD.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/961c924b108f55c424dcc16385e8aba5f1224256/e2e/$guid1.md", "", "", "$guid1.md")

$loD.ListRows.Add() | Out-Null
$wsD.Range("B3").Value = ".md"
$wsD.Range("C3").Value = "Handed back: in sync with en-US"
$wsD.Range("D3").Value = "e2e"
$wsD.Range("E3").Value = "ht"
$wsD.Range("F3").Value = "'True"
$wsD.Range("G3").Value = $dedeXlf2
$wsD.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("H3").Value = "2016-11-15 17:38:26"
$wsD.Range("J3").Value = $dedeXlf2
$wsD.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsD.Range("K3").Value = "2016-11-15 17:39:18"
$wsD.Range("L3").Value = "'"
$wsD.Range("M3").Value = "'True"
$wsD.Range("N3").Value = "'"
$wsD.Range("O3").Value = "'True"
$wsD.Range("P3").Value = "'"
$wsD.Range("A3").Value = "$guid2.md"
$wsD.Range("I3").Value = "$guid2.md"
$wsD.Hyperlinks.Add($wsD.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$guid2.md", "", "", "$guid2.md")
$wsD.Hyperlinks.Add($wsD.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/961c924b108f55c424dcc16385e8aba5f1224256/e2e/$guid2.md", "", "", "$guid2.md")

Write-Host "Report generated for handback."
